$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append starting at row 54 (A=test case, B=status, C=browser)
$rows = @(
    @("Upload PDF file to the system as attach", "PASSED", "chrome"),
    @("Student should see PDF file to the system as attach", "PASSED", "chrome"),
    @("Add new live session", "FAILED", "chrome"),
    @("Add new live session", "FAILED", "chrome"),
    @("Add new live session", "FAILED", "chrome"),
    @("Student should see the teacher's add a new live session", "FAILED", "chrome"),
    @("Taking an excel report of grades", "FAILED", "chrome"),
    @("Taking an excel report of grades", "PASSED", "chrome"),
    @("Add new live session", "FAILED", "chrome"),
    @("Student should see the teacher's add a new live session", "FAILED", "chrome"),
    @("Add new live session", "FAILED", "chrome"),
    @("Student should see the teacher's add a new live session", "FAILED", "chrome"),
    @("Add new live session", "FAILED", "chrome"),
    @("Student should see the teacher's add a new live session", "FAILED", "chrome"),
    @("Upload PDF file to the system as attach", "PASSED", "chrome"),
    @("Student should see PDF file to the system as attach", "PASSED", "chrome"),
    @("Add new live session", "PASSED", "chrome"),
    @("Student should see the teacher's add a new live session", "PASSED", "chrome"),
    @("Upload PDF file to the system as attach", "PASSED", "chrome"),
    @("Student should see PDF file to the system as attach", "PASSED", "chrome"),
    @("Add new live session", "PASSED", "chrome"),
    @("Student should see the teacher's add a new live session", "PASSED", "chrome")
)

$startRow = 54
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
}
